$wb = $excel.ActiveWorkbook

# --- mainoptions sheet: update pension/savings figures ---
$wsMain = $wb.Worksheets.Item("mainoptions")
$wsMain.Range("B1").Value = 1008
$wsMain.Range("B2").Value = 1368
$wsMain.Range("B3").Value = 1608
$wsMain.Range("B5").Value = 2484
$wsMain.Range("C9").Select()

# --- categories sheet: leave values as-is, selection stays put ---
$wsCategories = $wb.Worksheets.Item("categories")

# --- breakdowns sheet: update the pension/savings columns (I & J) ---
$wsBreakdowns = $wb.Worksheets.Item("breakdowns")
$wsBreakdowns.Range("I1").Value = 10
$wsBreakdowns.Range("J1").Value = 10
$wsBreakdowns.Range("I2").Value = 10
$wsBreakdowns.Range("J2").Value = 10
$wsBreakdowns.Range("I3").Value = 10
$wsBreakdowns.Range("J3").Value = 10
$wsBreakdowns.Range("I4").Value = 10
$wsBreakdowns.Range("J4").Value = 10
$wsBreakdowns.Range("I5").Value = 10
$wsBreakdowns.Range("J5").Value = 10
$wsBreakdowns.Range("I6").Value = 25
$wsBreakdowns.Range("J6").Value = 25

# Make "breakdowns" the active sheet, with its own selection
$wsBreakdowns.Activate()
$wsBreakdowns.Range("J7").Select()
